$p = $ppt.ActivePresentation
$s = $p.Slides.Item(18)

# "Picture 3" (shape id 1027) is enlarged/repositioned on the title slide
$pic = $s.Shapes.Item("Picture 3")
$pic.Left = 170.2394881889764
$pic.Top = 43.961692913385825
$pic.Width = 379.52114173228347
$pic.Height = 417.9162598425197
